$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -0.3332622562625545
$ws.Range("D2").Value = 0.7420916008870835

$ws.Range("C3").Value = -1.767954405700968
$ws.Range("D3").Value = 0.09093616888524236

$ws.Range("C4").Value = -0.5626075066492342
$ws.Range("D4").Value = 0.5793907117962123

$ws.Range("C5").Value = -1.962678472679351
$ws.Range("D5").Value = 0.06245213556464368

$ws.Range("C6").Value = -1.223936159792576
$ws.Range("D6").Value = 0.2339229688961586

$ws.Range("C7").Value = -0.3069587242118178
$ws.Range("D7").Value = 0.7617614259192389

$ws.Range("C8").Value = -1.167274145846048
$ws.Range("D8").Value = 0.2555960509960471

$ws.Range("C9").Value = 0.6883565558837057
$ws.Range("D9").Value = 0.4984232261285846

$ws.Range("C10").Value = 0.0949911637063853
$ws.Range("D10").Value = 0.9251817585945477

$ws.Range("C11").Value = -0.8141616532026933
$ws.Range("D11").Value = 0.4242777250696776
